# Weekly refresh of the fruit/vegetable (Berenjena) price records.
# The reported rows (2-25) get their per-record fields (Fecha, Calidad,
# Volumen, Precio minimo/maximo/promedio, Unidad de comercializacion,
# Precio $/Kg, Kg o Unidades) reshuffled across the existing rows, while
# the market/region/category columns (A, B, C, E, F, G, H, O, R) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as one "record" per row.
$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")

# target row -> source row (i.e. row $target ends up with the values that
# used to live in row $source)
$rowMap = @{
    2  = 6
    3  = 22
    4  = 10
    5  = 13
    6  = 16
    7  = 20
    8  = 14
    9  = 18
    10 = 17
    11 = 12
    12 = 24
    13 = 8
    14 = 25
    15 = 23
    16 = 2
    17 = 19
    18 = 5
    19 = 3
    20 = 15
    21 = 11
    22 = 9
    23 = 4
    24 = 21
    25 = 7
}

# Snapshot every source row's values before writing anything, so the
# shuffle is safe regardless of write order.
$snapshot = @{}
foreach ($row in $rowMap.Keys) {
    $values = @{}
    foreach ($col in $cols) {
        $values[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $values
}

foreach ($target in $rowMap.Keys) {
    $source = $rowMap[$target]
    $values = $snapshot[$source]
    foreach ($col in $cols) {
        $ws.Range("$col$target").Value = $values[$col]
    }
}
